# Add new rows 557-598 with NTP molecule data (TEST predictions)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @('CC(COc1ccccc1)OCCCO', 0.17, 0, '51730-94-0', 'Dipropylene glycol phenyl ether', 0),
    @('Cc1ccccc1C(=O)O', -0.05, 0, '118-90-1', 'o-Toluic acid', 0),
    @('Cc1cccc(C)c1O', 0.06, 0, '576-26-1', '2,6-Dimethyl phenol', 0),
    @('Cc1ccc2c(C)cccc2c1', 0.16, 0, '575-43-9', '1,6-Dimethylnaphthalene', 0),
    @('CC(C)CCOC(=O)/C=C/c1ccccc1', 0.08, 0, '7779-65-9', 'Isoamyl cinnamate', 0),
    @('Cc1ccc(C(=O)O)cc1', 0.34, 0, '99-94-5', 'p-Toluic acid', 0),
    @('OCC1OC(Oc2ccc(O)cc2)C(O)C(O)C1O', 0.26, 0, '497-76-7', 'Arbutin', 0),
    @('Cc1ccccc1O', 0.04, 0, '95-48-7', 'o-Cresol', 0),
    @('CC(C)c1ccccc1', 0.12, 0, '98-82-8', 'Cumene', 0),
    @('Cc1ccc(C(C)(C)O)cc1', 0.14, 0, '1197-01-9', 'p-Cymen-8-ol', 0),
    @('Cc1ccc2ccccc2c1C', 0.01, 0, '28804-88-8', 'Dimethylnaphthalene', 0),
    @('O=C(O)c1ccc(O)cc1O', 0.16, 0, '89-86-1', '2,4-Dihydroxybenzoic acid', 0),
    @('COc1cc(C2Oc3cc(C4Oc5cc(O)cc(O)c5C(=O)C4O)ccc3OC2CO)ccc1O', 0.41, 0, '22888-70-6', 'Silybin', 0),
    @('O=C(Cc1ccccc1)OCc1ccccc1', 0.04, 0, '102-16-9', 'Benzyl phenylacetate', 0),
    @('COCCOC(=O)c1ccccc1C(=O)OCCOC', 0.36, 1, '117-82-8', 'Di(2-methoxyethyl)phthalate', 0),
    @('CC(=O)c1cc(C)oc1C', 0.15, 1, '10599-70-9', '3-Acetyl-2,5-dimethylfuran', 0),
    @('Cc1cc(C)cc(O)c1', 0.07, 0, '108-68-9', '3,5-Dimethyl phenol', 0),
    @('Cc1ccc(O)c(O)c1', 0.41, 0, '452-86-8', 'p-Methylcatechol', 0),
    @('O=Cc1ccccc1C=O', 0.16, 1, '643-79-8', 'ortho-Phthalaldehyde', 0),
    @('CCCOC(=O)c1ccc(O)cc1', 0.05, 0, '94-13-3', 'Propyl-4-hydroxybenzoate', 0),
    @('Oc1ccc(/C=C/c2cc(O)cc(O)c2)cc1', 0.83, 0, '501-36-0', 'Resveratrol', 1),
    @('CC(CCC(=O)O)(c1ccc(O)cc1)c1ccc(O)cc1', -0.18, 0, '126-00-1', 'Diphenolic acid', 0),
    @('OCCOc1ccccc1', 0.08, 0, '122-99-6', 'Ethylene glycol monophenyl ether', 0),
    @('Cc1cccc(C(=O)O)c1', 0.06, 0, '99-04-7', 'm-Toluic acid', 0),
    @('CCCCOC(=O)c1ccc(O)cc1', -0.05, 0, '94-26-8', 'Butyl Paraben', 0),
    @('Cc1ccc(O)cc1C', 0.31, 0, '95-65-8', '3,4-Dimethyl phenol', 0),
    @('Cc1cccc2ccccc12', -0.02, 0, '1321-94-4', 'Methyl naphthalene', 0),
    @('Cc1ccc(C)c(O)c1', 0.25, 0, '95-87-4', '2,5-Dimethyl phenol', 0),
    @('Cc1ccc(O)c(C)c1', 0.47, 1, '105-67-9', '2,4-Dimethylphenol', 0),
    @('CC(C)OC(=O)Cc1ccccc1', 0.03, 0, '4861-85-2', 'Isopropyl phenylacetate', 0),
    @('C=CCc1cc(OC)c2c(c1)OCO2', 0.07, 0, '607-91-0', 'Myristicin', 0),
    @('O=C(O)/C=C/c1ccc(/C=C/C(=O)O)cc1', 0.41, 0, '16323-43-6', '3,3''-(1,4-Phenylene)bis-2-propenoic acid', 0),
    @('CC(O)COc1ccccc1', 0.21, 0, '770-35-4', 'Propylene glycol phenyl ether', 0),
    @('CCCCCCOC(=O)c1ccccc1', -0.03, 0, '6789-88-4', 'Hexylbenzoate', 0),
    @('Cc1c(O)cccc1O', 0.16, 0, '608-25-3', '2-Methyl-1,3-benzenediol', 0),
    @('Cc1ccc2ccccc2c1', 0.04, 0, '91-57-6', '2-Methylnaphthalene', 0),
    @('CC(=O)OC(C(=O)c1ccccc1)c1ccccc1', 0.2, 0, '574-06-1', 'Benzoin acetate', 0),
    @('C=Cc1cccc(C)c1', 0.01, 0, '25013-15-4', 'Vinyl toluene', 0),
    @('Cc1ccc2cc3c(ccc4ccccc43)c3c2c1CC3', 0.91, 1, '56-49-5', '3-Methylcholanthrene', 1),
    @('O=c1cc(-c2ccccc2)oc2ccc3ccccc3c12', 0.13, 0, '6051-87-2', '5,6-Benzoflavone', 0),
    @('O=c1c(-c2ccc(O)cc2)coc2cc(O)cc(O)c12', 0.17, 1, '446-72-0', 'Genistein', 0),
    @('Cc1ccc2oc(=O)ccc2c1', 0.38, 0, '92-48-8', 'Methyl coumarin', 0),
)

$startRow = 557
for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $row = $newRows[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
}

# Column width adjustments (character-width units; Excel stores width
# with an internal padding offset of 5/6 relative to the ColumnWidth API)
$ws.Columns.Item(2).ColumnWidth = (26 - 5/6)
$ws.Columns.Item(6).ColumnWidth = (15 - 5/6)
